$d = $word.ActiveDocument

# Curly double quotes used around "Partes"
$lq = [char]0x201C
$rq = [char]0x201D

# --- Pass 1: append all paragraphs with their final text, using the
# default (Normal) paragraph style. Styling is applied afterwards in a
# second pass so that body paragraphs never inherit a heading style. ---

$p1 = $d.Paragraphs.Add()
$p1.Range.Text = "ACUERDO DE CONFIDENCIALIDAD (NDA)"

$p2 = $d.Paragraphs.Add()
$p2.Range.Text = "Entre FUNDADOR 1, FUNDADOR 2 y cualquier tercero que se adhiera (las " + $lq + "Partes" + $rq + "), con fecha ___ de __________ de 2025."

$p3 = $d.Paragraphs.Add()
$p3.Range.Text = "1. Información Confidencial"

$p4 = $d.Paragraphs.Add()
$p4.Range.Text = "Incluye código, prototipos, planes de negocio y datos de clientes."

$p5 = $d.Paragraphs.Add()
$p5.Range.Text = "2. Obligaciones"

$p6 = $d.Paragraphs.Add()
$p6.Range.Text = "La información se usará sólo para la cooperación y no se divulgará sin permiso escrito."

$p7 = $d.Paragraphs.Add()
$p7.Range.Text = "3. Excepciones"

$p8 = $d.Paragraphs.Add()
$p8.Range.Text = "No confidencial: dominio público, posesión previa legítima o requerimiento judicial."

$p9 = $d.Paragraphs.Add()
$p9.Range.Text = "4. Vigencia"

$p10 = $d.Paragraphs.Add()
$p10.Range.Text = "Obligaciones vigentes por 5 años desde la última divulgación."

$p11 = $d.Paragraphs.Add()
$p11.Range.Text = "5. Ley y Arbitraje"

$p12 = $d.Paragraphs.Add()
$p12.Range.Text = "Ley argentina. Arbitraje en CABA."

# --- Pass 2: apply heading styles to the title + section-heading
# paragraphs (1, 3, 5, 7, 9, 11). The remaining body paragraphs keep the
# default Normal style. ---

$d.Paragraphs(1).Style = "Heading1"
$d.Paragraphs(3).Style = "Heading2"
$d.Paragraphs(5).Style = "Heading2"
$d.Paragraphs(7).Style = "Heading2"
$d.Paragraphs(9).Style = "Heading2"
$d.Paragraphs(11).Style = "Heading2"
